$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the migration wave date in the title cell
$ws.Range("A1").Value = "Hotcarding Spreadsheet - Migration Wave 10/24/2002"

# Append a new data row (row 5) under the existing data row (row 4)
$ws.Range("B5").Value = "NewFinance Ltd"
$ws.Range("C5").Value = "ENT99999"
$ws.Range("D5").Value = "FISN"
$ws.Range("E5").Value = "OldBank"
$ws.Range("F5").Value = "PaymentsNext"
$ws.Range("G5").Value = "Standard"
$ws.Range("H5").Value = "Bangalore"
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = "Yes"

# A5 looks like a date ("2025-11-06"); force it to remain plain text (matches
# row 4's A4, which is also stored as literal text, not a date serial).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-11-06"
# Re-align A5's cell format with the rest of the row (avoid leaving a stray
# "text" number-format style behind once the literal-text value is locked in).
$ws.Range("B5").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
